$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows at 843..848, shifting existing rows 843-913 down to 849-919.
$ws.Range("A843:A848").EntireRow.Insert()

# Constant column values shared by all rows in this block.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102005
$categoria = "Naranja"
$unidad = "`$/malla 13 kilos"
$origen = "Provincia de Quillota"
$kgUnidad = 13

# New rows data: row -> (Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg)
$newRows = @(
    @{ Row = 843; Variedad = "Fukumoto";   Calidad = "Primera"; Volumen = 172; PMin = 3000; PMax = 3500; PProm = 3247; PKg = 250 },
    @{ Row = 844; Variedad = "Fukumoto";   Calidad = "Segunda"; Volumen = 145; PMin = 2400; PMax = 2500; PProm = 2455; PKg = 189 },
    @{ Row = 845; Variedad = "Navel Late"; Calidad = "Primera"; Volumen = 142; PMin = 3000; PMax = 3500; PProm = 3264; PKg = 251 },
    @{ Row = 846; Variedad = "Navel Late"; Calidad = "Segunda"; Volumen = 75;  PMin = 2500; PMax = 2500; PProm = 2500; PKg = 192 },
    @{ Row = 847; Variedad = "New Hall";   Calidad = "Primera"; Volumen = 214; PMin = 3000; PMax = 3500; PProm = 3236; PKg = 249 },
    @{ Row = 848; Variedad = "New Hall";   Calidad = "Segunda"; Volumen = 138; PMin = 2400; PMax = 2500; PProm = 2449; PKg = 188 }
)

$fechaNueva = 44783

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fechaNueva
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
